$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1722.50748368523
$ws.Range("C2").Value = 1128.65460342033
$ws.Range("D2").Value = 814.287878069938
$ws.Range("E2").Value = 2316.36036395013
$ws.Range("F2").Value = 2630.72708930052
$ws.Range("I2").Value = 523.50748368523

$ws.Range("B3").Value = 1573.7547520467
$ws.Range("C3").Value = 792.728171138095
$ws.Range("D3").Value = 379.277673866225
$ws.Range("E3").Value = 2354.78133295532
$ws.Range("F3").Value = 2768.23183022718
$ws.Range("I3").Value = 1267.7547520467

$ws.Range("B4").Value = 5000.35717921842
$ws.Range("C4").Value = 1905.97617927905
$ws.Range("D4").Value = 267.909826037655
$ws.Range("E4").Value = 8094.73817915779
$ws.Range("F4").Value = 9732.80453239918
$ws.Range("I4").Value = 4803.35717921842

$ws.Range("B5").Value = 8709.16575577456
$ws.Range("C5").Value = 2371.80108461244
$ws.Range("D5").Value = -982.997124223106
$ws.Range("E5").Value = 15046.5304269367
$ws.Range("F5").Value = 18401.3286357722
$ws.Range("I5").Value = 8593.16575577456

$ws.Range("B6").Value = 9717.02457056929
$ws.Range("C6").Value = 1665.09733903548
$ws.Range("D6").Value = -2597.33546318961
$ws.Range("E6").Value = 17768.9518021031
$ws.Range("F6").Value = 22031.3846043282
$ws.Range("I6").Value = 9627.02457056929

$ws.Range("B7").Value = 5505.56319948048
$ws.Range("C7").Value = 414.285516131601
$ws.Range("D7").Value = -2280.87408711893
$ws.Range("E7").Value = 10596.8408828294
$ws.Range("F7").Value = 13292.0004860799
$ws.Range("I7").Value = 5391.56319948048

$ws.Range("B8").Value = 2489.08076783683
$ws.Range("C8").Value = -44.2577403004939
$ws.Range("D8").Value = -1385.32613867542
$ws.Range("E8").Value = 5022.41927597416
$ws.Range("F8").Value = 6363.48767434909
$ws.Range("I8").Value = 2393.08076783683

$ws.Range("B9").Value = 1236.72874640402
$ws.Range("C9").Value = -134.68338178881
$ws.Range("D9").Value = -860.665109205004
$ws.Range("E9").Value = 2608.14087459684
$ws.Range("F9").Value = 3334.12260201304
$ws.Range("I9").Value = 1154.72874640402

$ws.Range("B10").Value = 1020.38472702269
$ws.Range("C10").Value = -203.003155959862
$ws.Range("D10").Value = -850.625581290224
$ws.Range("E10").Value = 2243.77261000524
$ws.Range("F10").Value = 2891.3950353356
$ws.Range("I10").Value = 918.384727022687

$ws.Range("B11").Value = 1086.38153310003
$ws.Range("C11").Value = -313.458742176718
$ws.Range("D11").Value = -1054.4894216205
$ws.Range("E11").Value = 2486.22180837678
$ws.Range("F11").Value = 3227.25248782056
$ws.Range("I11").Value = 929.38153310003

$ws.Range("B12").Value = 1268.69925915799
$ws.Range("C12").Value = -479.756132394754
$ws.Range("D12").Value = -1405.33250683028
$ws.Range("E12").Value = 3017.15465071074
$ws.Range("F12").Value = 3942.73102514626
$ws.Range("I12").Value = 992.699259157991

$ws.Range("B13").Value = 1718.01420378355
$ws.Range("C13").Value = -804.340402051133
$ws.Range("D13").Value = -2139.59427377611
$ws.Range("E13").Value = 4240.36880961823
$ws.Range("F13").Value = 5575.62268134321
$ws.Range("I13").Value = 1614.01420378355

$ws.Range("B14").Value = 1722.52003769654
$ws.Range("C14").Value = -962.848490676093
$ws.Range("D14").Value = -2384.39672034359
$ws.Range("E14").Value = 4407.88856606918
$ws.Range("F14").Value = 5829.43679573668
$ws.Range("I14").Value = 1595.52003769654

$ws.Range("B15").Value = 1573.76622191557
$ws.Range("C15").Value = -1024.18726743177
$ws.Range("D15").Value = -2399.4607706847
$ws.Range("E15").Value = 4171.71971126291
$ws.Range("F15").Value = 5546.99321451584
$ws.Range("I15").Value = 1527.76622191557

$ws.Range("B16").Value = 5000.39362291529
$ws.Range("C16").Value = -3719.68604150992
$ws.Range("D16").Value = -8335.81737639056
$ws.Range("E16").Value = 13720.4732873405
$ws.Range("F16").Value = 18336.6046222211
$ws.Range("I16").Value = 4918.39362291529

$ws.Range("B17").Value = 8709.22923007958
$ws.Range("C17").Value = -7302.47465403387
$ws.Range("D17").Value = -15778.5586886581
$ws.Range("E17").Value = 24720.933114193
$ws.Range("F17").Value = 33197.0171488173
$ws.Range("I17").Value = 8602.22923007958

$ws.Range("B18").Value = 9717.0953903698
$ws.Range("C18").Value = -9083.3900349258
$ws.Range("D18").Value = -19035.7658413909
$ws.Range("E18").Value = 28517.5808156654
$ws.Range("F18").Value = 38469.9566221305
$ws.Range("I18").Value = 9530.0953903698

$ws.Range("B19").Value = 5505.60332522934
$ws.Range("C19").Value = -5687.28189116755
$ws.Range("D19").Value = -11612.4374130738
$ws.Range("E19").Value = 16698.4885416262
$ws.Range("F19").Value = 22623.6440635325
$ws.Range("I19").Value = 5390.60332522934

$ws.Range("B20").Value = 2489.09890880192
$ws.Range("C20").Value = -2820.89896545507
$ws.Range("D20").Value = -5631.84204382005
$ws.Range("E20").Value = 7799.09678305891
$ws.Range("F20").Value = 10610.0398614239
$ws.Range("I20").Value = 2387.09890880192

$ws.Range("B21").Value = 1236.73775995364
$ws.Range("C21").Value = -1528.44338916086
$ws.Range("D21").Value = -2992.241868104
$ws.Range("E21").Value = 4001.91890906813
$ws.Range("F21").Value = 5465.71738801127
$ws.Range("I21").Value = 1168.73775995364

$ws.Range("B22").Value = 1020.39216380977
$ws.Range("C22").Value = -1368.20900618701
$ws.Range("D22").Value = -2632.65808597434
$ws.Range("E22").Value = 3408.99333380655
$ws.Range("F22").Value = 4673.44241359388
$ws.Range("I22").Value = 930.39216380977

$ws.Range("B23").Value = 1086.38945088627
$ws.Range("C23").Value = -1573.59590189417
$ws.Range("D23").Value = -2981.70709007363
$ws.Range("E23").Value = 3746.37480366671
$ws.Range("F23").Value = 5154.48599184617
$ws.Range("I23").Value = 985.389450886272

$ws.Range("B24").Value = 1268.7085057157
$ws.Range("C24").Value = -1977.69438744411
$ws.Range("D24").Value = -3696.23626979628
$ws.Range("E24").Value = 4515.11139887551
$ws.Range("F24").Value = 6233.65328122768
$ws.Range("I24").Value = 1149.7085057157

$ws.Range("B25").Value = 1718.02672504686
$ws.Range("C25").Value = -2872.73189195016
$ws.Range("D25").Value = -5302.93270642748
$ws.Range("E25").Value = 6308.78534204388
$ws.Range("F25").Value = 8738.98615652119
$ws.Range("I25").Value = 1560.02672504686

$ws.Range("B26").Value = 1722.53259179935
$ws.Range("C26").Value = -3080.76099981079
$ws.Range("D26").Value = -5623.47103367886
$ws.Range("E26").Value = 6525.82618340949
$ws.Range("F26").Value = 9068.53621727756
$ws.Range("I26").Value = 1652.53259179935

$ws.Range("B27").Value = 1573.77769186803
$ws.Range("C27").Value = -3002.97512390532
$ws.Range("D27").Value = -5425.76171507428
$ws.Range("E27").Value = 6150.53050764139
$ws.Range("F27").Value = 8573.31709881035
$ws.Range("I27").Value = 1518.77769186803

$ws.Range("B28").Value = 5000.43006687777
$ws.Range("C28").Value = -10156.7293834669
$ws.Range("D28").Value = -18180.4449203742
$ws.Range("E28").Value = 20157.5895172224
$ws.Range("F28").Value = 28181.3050541297
$ws.Range("I28").Value = 4902.43006687777

$ws.Range("B29").Value = 8709.29270484721
$ws.Range("C29").Value = -18792.8318357884
$ws.Range("D29").Value = -33351.5771401165
$ws.Range("E29").Value = 36211.4172454828
$ws.Range("F29").Value = 50770.1625498109
$ws.Range("I29").Value = 8612.29270484721

$ws.Range("B30").Value = 9717.16621068647
$ws.Range("C30").Value = -22234.4197383628
$ws.Range("D30").Value = -39148.5676331642
$ws.Range("E30").Value = 41668.7521597358
$ws.Range("F30").Value = 58582.9000545372
$ws.Range("I30").Value = 9562.16621068647

$ws.Range("B31").Value = 5505.64345127063
$ws.Range("C31").Value = -13337.1253935064
$ws.Range("D31").Value = -23311.8846901603
$ws.Range("E31").Value = 24348.4122960477
$ws.Range("F31").Value = 34323.1715927015
$ws.Range("I31").Value = 5338.64345127063

$ws.Range("B32").Value = 2489.11704989922
$ws.Range("C32").Value = -6374.19923132671
$ws.Range("D32").Value = -11066.1554509058
$ws.Range("E32").Value = 11352.4333311251
$ws.Range("F32").Value = 16044.3895507042
$ws.Range("I32").Value = 2383.11704989922

$ws.Range("B33").Value = 1236.74677356895
$ws.Range("C33").Value = -3343.52069589339
$ws.Range("D33").Value = -5768.16783233824
$ws.Range("E33").Value = 5817.01424303128
$ws.Range("F33").Value = 8241.66137947613
$ws.Range("I33").Value = 1108.74677356895

$ws.Range("B34").Value = 1020.39960065105
$ws.Range("C34").Value = -2908.7382576287
$ws.Range("D34").Value = -4988.69819889679
$ws.Range("E34").Value = 4949.53745893081
$ws.Range("F34").Value = 7029.4974001989
$ws.Range("I34").Value = 922.399600651055

$ws.Range("B35").Value = 1086.39736873022
$ws.Range("C35").Value = -3261.72889103639
$ws.Range("D35").Value = -5563.48789418066
$ws.Range("E35").Value = 5434.52362849683
$ws.Range("F35").Value = 7736.2826316411
$ws.Range("I35").Value = 968.39736873022

$ws.Range("B36").Value = 1268.7177523408
$ws.Range("C36").Value = -4007.77128069096
$ws.Range("D36").Value = -6800.97585001383
$ws.Range("E36").Value = 6545.20678537256
$ws.Range("F36").Value = 9338.41135469543
$ws.Range("I36").Value = 1148.7177523408
